$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Oklo Inc.) updates
$ws.Range("K2").Value = 58.5
$ws.Range("N2").Value = 51.53902399942638

# Row 3 (NuScale Power Corporation) updates
$ws.Range("K3").Value = 52.9
$ws.Range("N3").Value = 51.53902399942638
